# Add code for discretized normal distribution:
# introduces a new "gen_html_10x10" sheet, adds a "Field 12" column and a
# second header row (pairnum/halfpos/mat_1_.../mat_2_.../treat/tasknum) to
# the "10x10" sheet, and nudges the active selections.

$wb = $excel.ActiveWorkbook

$dict = $wb.Worksheets.Item("dict")
$grid = $wb.Worksheets.Item("10x10")

# --- "10x10" sheet: extend header row and add a second header row ---
$grid.Range("L1").Value = "Field 12"

$grid.Range("A2").Value = "pairnum"
$grid.Range("B2").Value = "halfpos"
$grid.Range("C2").Value = "mat_1_half_tag"
$grid.Range("D2").Value = "mat_1_whole_tag"
$grid.Range("E2").Value = "mat_1_half_red"
$grid.Range("F2").Value = "mat_1_total_red"
$grid.Range("G2").Value = "mat_2_half_tag"
$grid.Range("H2").Value = "mat_2_whole_tag"
$grid.Range("I2").Value = "mat_2_half_red"
$grid.Range("J2").Value = "mat_2_total_red"
$grid.Range("K2").Value = "treat"
$grid.Range("L2").Value = "tasknum"

# --- add the new "gen_html_10x10" worksheet after "10x10" ---
$newSheet = $wb.Worksheets.Add($null, $grid)
$newSheet.Name = "gen_html_10x10"

$newSheet.Columns.Item(1).ColumnWidth = 75.33333333333334
$newSheet.Columns.Item(2).ColumnWidth = 74.5
$newSheet.Columns.Item(3).ColumnWidth = 12.333333333333332
$newSheet.Columns.Item(4).ColumnWidth = 14.166666666666666

$newSheet.Activate() | Out-Null
$newSheet.Range("B8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 168

# --- restore selections / active sheet to match the target state ---
$dict.Activate() | Out-Null
$dict.Range("D3").Select() | Out-Null

$grid.Activate() | Out-Null
$grid.Range("C17").Select() | Out-Null
